# Performance file by manik on 16.11.2018.15.54
#
# Fills in the Mobile / Desktop PageSpeed-Insight scores on the
# "Page_Insight" sheet (columns B and C, rows 2-6), which previously held
# empty, but pre-formatted, cells.
#
# The target cells are written with Excel's ordinary text-entry semantics
# (stored as shared-string text, like `Urls`, `Mobile`, etc.), not as
# numbers -- even though a couple of rows carry a percentage number
# format. Assigning `.Value`/`.Value2` directly lets Excel's COM layer
# auto-coerce the numeric-looking text into a real number (and, if we try
# to force text via NumberFormat = "@", Excel mints a brand-new cell
# style instead of reusing the existing one). To avoid both side effects
# we stage the text as a literal-string formula and then collapse it to a
# plain value in place via Copy + PasteSpecial(xlPasteValues), which keeps
# the original cell style untouched while leaving a plain text value
# behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Page_Insight")

$xlPasteValues = -4163

$scores = [ordered]@{
    "B2" = "47";  "C2" = "98";
    "B3" = "84";  "C3" = "100";
    "B4" = "84";  "C4" = "100";
    "B5" = "26";  "C5" = "82";
    "B6" = "14";  "C6" = "65";
}

foreach ($addr in $scores.Keys) {
    $cell = $ws.Range($addr)
    $cell.Formula = '="' + $scores[$addr] + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial($xlPasteValues) | Out-Null
}

$excel.CutCopyMode = 0
